$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.748.68"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "3.489.94"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "606.05"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").Value = "192.40"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").Value = "53.33"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("D13").Value = "9.58"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "4.054.65"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").Value = "610.95"
$ws.Range("E15").Value = "  +4.99%  "
$ws.Range("D16").Value = "69.863.87"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "12.65"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "18.81"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "3.509.36"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "0.988"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "17.78"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "105.09"
$ws.Range("E23").Value = "  +10.93%  "
$ws.Range("D24").Value = "4.63"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("E26").Value = "  +4.12%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  +5.40%  "
$ws.Range("D29").Value = "34.20"
$ws.Range("E29").Value = "  +5.67%  "
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "4.27"
$ws.Range("E31").Value = "  +13.40%  "
$ws.Range("D32").Value = "12.64"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "64.20"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "3.742.16"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "3.06"
$ws.Range("E37").Value = "  -5.10%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "517.94"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("D39").Value = "0.0₃0791"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").Value = "36.54"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("D46").Value = "0.140"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").Value = "8.73"
$ws.Range("E49").Value = "  -5.34%  "
$ws.Range("D50").Value = "132.68"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "1.29"
$ws.Range("E51").Value = "  +10.65%  "
